# Updated cryptos list (price & 1h-volume columns) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Subscript-six glyph used in BabyDogeCoin's price (e.g. 0.0<sub>6</sub>0251).
$sub6 = [char]0x2086

$ws.Range("D2").Value = '65.524.55'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '2.646.35'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.84'
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.53'
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.629'
$ws.Range("E8").Value = '  +4.04%  '
$ws.Range("E9").Value = '  +3.10%  '
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  -3.25%  '
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.65'
$ws.Range("E13").Value = '  -2.87%  '
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D15").Value = '3.119.86'
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = '65.387.73'
$ws.Range("D17").Value = '2.622.04'
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.56'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.42'
$ws.Range("E20").Value = '  -2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '347.79'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.90'
$ws.Range("E23").Value = '  -2.39%  '
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.57'
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("E27").Value = '  -3.05%  '
$ws.Range("E28").Value = '  -3.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.01'
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '538.56'
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.87'
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.09'
$ws.Range("E32").Value = '  -3.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.74'
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("E36").Value = '  -1.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.24'
$ws.Range("E37").Value = '  -1.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '154.73'
$ws.Range("E39").Value = '  -3.60%  '
$ws.Range("E40").Value = '  -2.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '160.24'
$ws.Range("E42").Value = '  -3.63%  '
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("E44").Value = '  +2.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0601'
$ws.Range("E45").Value = '  -2.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.45'
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.633'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0254'
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0999'
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").Value = "0.0${sub6}0251"
$ws.Range("E50").Value = '  +6.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.56'
$ws.Range("E51").Value = '  -3.96%  '
